$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# COM type inference are pre-formatted as Text so the stored value stays a string.
$ws.Range("D2").Value = "57.441.44"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.012.83"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.39"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.57"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "3.528.57"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.36"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "57.447.46"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.22"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "3.008.75"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.11"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.499"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.38"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.77"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.33"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.56"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.85"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.63"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0677"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").Value = "3.046.78"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("E42").Value = "  +4.68%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "2.223.69"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.05"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.54"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  -5.84%  "
